$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 34.71251733333333
$ws.Range("H2").Value = 104.137552
$ws.Range("I2").Value = 0.111750244749681
$ws.Range("J2").Value = 0.111750244749681
$ws.Range("M2").Value = 2.565830333333333
$ws.Range("N2").Value = 7.697490999999999
$ws.Range("O2").Value = 0.0934185609347503
$ws.Range("P2").Value = 0.0934185609347503
$ws.Range("Q2").Value = 89.06642992022576
$ws.Range("R2").Value = 801.597869282032
$ws.Range("S2").Value = 0.01043954704862134
$ws.Range("T2").Value = 0.01043954704862134

$ws.Range("G3").Value = 34.71251733333333
$ws.Range("H3").Value = 104.137552
$ws.Range("I3").Value = 0.111750244749681
$ws.Range("J3").Value = 0.111750244749681
$ws.Range("O3").Value = 0.3847798091300315
$ws.Range("P3").Value = 0.3847798091300315
$ws.Range("Q3").Value = 366.8539052805031
$ws.Range("R3").Value = 3301.685147524528
$ws.Range("S3").Value = 0.04299923784501657
$ws.Range("T3").Value = 0.04299923784501657

$ws.Range("G4").Value = 34.71251733333333
$ws.Range("H4").Value = 104.137552
$ws.Range("I4").Value = 0.111750244749681
$ws.Range("J4").Value = 0.111750244749681
$ws.Range("M4").Value = 13.68376133333333
$ws.Range("N4").Value = 41.051284
$ws.Range("O4").Value = 0.4982080363333638
$ws.Range("P4").Value = 0.4982080363333638
$ws.Range("Q4").Value = 474.9978024685298
$ws.Range("R4").Value = 4274.980222216768
$ws.Range("S4").Value = 0.05567486999651138
$ws.Range("T4").Value = 0.05567486999651138

$ws.Range("G5").Value = 34.71251733333333
$ws.Range("H5").Value = 104.137552
$ws.Range("I5").Value = 0.111750244749681
$ws.Range("J5").Value = 0.111750244749681
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6480206666666667
$ws.Range("N5").Value = 1.944062
$ws.Range("O5").Value = 0.0235935936018545
$ws.Range("P5").Value = 0.0235935936018545
$ws.Range("Q5").Value = 22.49442862402489
$ws.Range("R5").Value = 202.449857616224
$ws.Range("S5").Value = 0.002636589859531748
$ws.Range("T5").Value = 0.002636589859531748

$ws.Range("H6").Value = 578.4917909999999
$ws.Range("I6").Value = 0.620780861354714
$ws.Range("J6").Value = 0.6207808613547139
$ws.Range("M6").Value = 2.565830333333333
$ws.Range("N6").Value = 7.697490999999999
$ws.Range("O6").Value = 0.0934185609347503
$ws.Range("P6").Value = 0.0934185609347503
$ws.Range("Q6").Value = 494.7705949773755
$ws.Range("R6").Value = 4452.93535479638
$ws.Range("S6").Value = 0.05799245472359213
$ws.Range("T6").Value = 0.05799245472359212

$ws.Range("H7").Value = 578.4917909999999
$ws.Range("I7").Value = 0.620780861354714
$ws.Range("J7").Value = 0.6207808613547139
$ws.Range("O7").Value = 0.3847798091300315
$ws.Range("P7").Value = 0.3847798091300315
$ws.Range("Q7").Value = 2037.900532759427
$ws.Range("S7").Value = 0.2388639413436434
$ws.Range("T7").Value = 0.2388639413436433

$ws.Range("H8").Value = 578.4917909999999
$ws.Range("I8").Value = 0.620780861354714
$ws.Range("J8").Value = 0.6207808613547139
$ws.Range("M8").Value = 13.68376133333333
$ws.Range("N8").Value = 41.051284
$ws.Range("O8").Value = 0.4982080363333638
$ws.Range("P8").Value = 0.4982080363333638
$ws.Range("Q8").Value = 2638.647867112183
$ws.Range("R8").Value = 23747.83080400964
$ws.Range("S8").Value = 0.3092780139288662
$ws.Range("T8").Value = 0.3092780139288662

$ws.Range("H9").Value = 578.4917909999999
$ws.Range("I9").Value = 0.620780861354714
$ws.Range("J9").Value = 0.6207808613547139
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6480206666666667
$ws.Range("N9").Value = 1.944062
$ws.Range("O9").Value = 0.0235935936018545
$ws.Range("P9").Value = 0.0235935936018545
$ws.Range("Q9").Value = 124.9582120216713
$ws.Range("R9").Value = 1124.623908195042
$ws.Range("S9").Value = 0.0146464513586123
$ws.Range("T9").Value = 0.0146464513586123

$ws.Range("G10").Value = 19.96051866666667
$ws.Range("H10").Value = 59.881556
$ws.Range("I10").Value = 0.06425903442584988
$ws.Range("J10").Value = 0.06425903442584986
$ws.Range("M10").Value = 2.565830333333333
$ws.Range("N10").Value = 7.697490999999999
$ws.Range("O10").Value = 0.0934185609347503
$ws.Range("P10").Value = 0.0934185609347503
$ws.Range("Q10").Value = 51.21530426399956
$ws.Range("R10").Value = 460.937738375996
$ws.Range("S10").Value = 0.006002986523119474
$ws.Range("T10").Value = 0.006002986523119473

$ws.Range("G11").Value = 19.96051866666667
$ws.Range("H11").Value = 59.881556
$ws.Range("I11").Value = 0.06425903442584988
$ws.Range("J11").Value = 0.06425903442584986
$ws.Range("O11").Value = 0.3847798091300315
$ws.Range("P11").Value = 0.3847798091300315
$ws.Range("Q11").Value = 210.9496742622982
$ws.Range("R11").Value = 1898.547068360684
$ws.Range("S11").Value = 0.02472557900125864
$ws.Range("T11").Value = 0.02472557900125863

$ws.Range("G12").Value = 19.96051866666667
$ws.Range("H12").Value = 59.881556
$ws.Range("I12").Value = 0.06425903442584988
$ws.Range("J12").Value = 0.06425903442584986
$ws.Range("M12").Value = 13.68376133333333
$ws.Range("N12").Value = 41.051284
$ws.Range("O12").Value = 0.4982080363333638
$ws.Range("P12").Value = 0.4982080363333638
$ws.Range("Q12").Value = 273.1349735242116
$ws.Range("R12").Value = 2458.214761717904
$ws.Range("S12").Value = 0.03201436735798069
$ws.Range("T12").Value = 0.03201436735798068

$ws.Range("G13").Value = 19.96051866666667
$ws.Range("H13").Value = 59.881556
$ws.Range("I13").Value = 0.06425903442584988
$ws.Range("J13").Value = 0.06425903442584986
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.6480206666666667
$ws.Range("N13").Value = 1.944062
$ws.Range("O13").Value = 0.0235935936018545
$ws.Range("P13").Value = 0.0235935936018545
$ws.Range("Q13").Value = 12.93482861338578
$ws.Range("R13").Value = 116.413457520472
$ws.Range("S13").Value = 0.001516101543491079
$ws.Range("T13").Value = 0.001516101543491079

$ws.Range("G14").Value = 63.12224
$ws.Range("H14").Value = 189.36672
$ws.Range("I14").Value = 0.2032098594697551
$ws.Range("J14").Value = 0.2032098594697551
$ws.Range("M14").Value = 2.565830333333333
$ws.Range("N14").Value = 7.697490999999999
$ws.Range("O14").Value = 0.0934185609347503
$ws.Range("P14").Value = 0.0934185609347503
$ws.Range("Q14").Value = 161.9609580999467
$ws.Range("R14").Value = 1457.64862289952
$ws.Range("S14").Value = 0.01898357263941737
$ws.Range("T14").Value = 0.01898357263941736

$ws.Range("G15").Value = 63.12224
$ws.Range("H15").Value = 189.36672
$ws.Range("I15").Value = 0.2032098594697551
$ws.Range("J15").Value = 0.2032098594697551
$ws.Range("O15").Value = 0.3847798091300315
$ws.Range("P15").Value = 0.3847798091300315
$ws.Range("Q15").Value = 667.0976936557865
$ws.Range("R15").Value = 6003.87924290208
$ws.Range("S15").Value = 0.0781910509401129
$ws.Range("T15").Value = 0.07819105094011289

$ws.Range("G16").Value = 63.12224
$ws.Range("H16").Value = 189.36672
$ws.Range("I16").Value = 0.2032098594697551
$ws.Range("J16").Value = 0.2032098594697551
$ws.Range("M16").Value = 13.68376133333333
$ws.Range("N16").Value = 41.051284
$ws.Range("O16").Value = 0.4982080363333638
$ws.Range("P16").Value = 0.4982080363333638
$ws.Range("Q16").Value = 863.7496669853867
$ws.Range("R16").Value = 7773.74700286848
$ws.Range("S16").Value = 0.1012407850500055
$ws.Range("T16").Value = 0.1012407850500055

$ws.Range("G17").Value = 63.12224
$ws.Range("H17").Value = 189.36672
$ws.Range("I17").Value = 0.2032098594697551
$ws.Range("J17").Value = 0.2032098594697551
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.6480206666666667
$ws.Range("N17").Value = 1.944062
$ws.Range("O17").Value = 0.0235935936018545
$ws.Range("P17").Value = 0.0235935936018545
$ws.Range("Q17").Value = 40.90451604629333
$ws.Range("R17").Value = 368.14064441664
$ws.Range("S17").Value = 0.004794450840219366
$ws.Range("T17").Value = 0.004794450840219366
